$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank columns before the existing column E ("AUTO.ROLLOVER").
# This shifts the old column E (with its header/value/width/bestFit) to column M,
# matching the target layout where AUTO.ROLLOVER now lives in column M.
$ws.Range("E1:L1").EntireColumn.Insert()

# Fill in the new header row cells (E1:L1, N1:O1) with the new field names.
$ws.Range("E1").Value = "INTEREST.RATE"
$ws.Range("F1").Value = "INTEND.DATE"
$ws.Range("G1").Value = "CUST.REMARKS:1"
$ws.Range("H1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("I1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("J1").Value = "PRIN.LIQ.ACCT"
$ws.Range("K1").Value = "INT.LIQ.ACCT"
$ws.Range("L1").Value = "CHRG.LIQ.ACCT"
$ws.Range("N1").Value = "FINAL.MATURITY"
$ws.Range("O1").Value = "EXP.DATE"

# Give the new columns E..L a uniform width (matches the other non-bestFit columns).
$ws.Range("E1:L1").ColumnWidth = 9.42

# Page setup: portrait orientation.
$ws.PageSetup.Orientation = 1

# Selection / active cell as left by the editing session.
$ws.Range("G8").Select()
